$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2 through 13 from serial 45175 to 45183
for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value2 = 45183
    }
}
